$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record needs to be inserted as row 140 (pushing the
# existing rows 140-146 down to 141-147). Inserting the row this way also
# carries over the date number-format style from the row above.
$ws.Rows.Item(140).Insert()

# Populate the newly inserted row 140 with the new "Zapallo italiano" record.
$ws.Cells.Item(140, 1).Value = 11
$ws.Cells.Item(140, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(140, 3).Value = "Bíobío"
$ws.Cells.Item(140, 4).Value = 44782
$ws.Cells.Item(140, 5).Value = 8
$ws.Cells.Item(140, 6).Value = 100112032
$ws.Cells.Item(140, 7).Value = "Zapallo italiano"
$ws.Cells.Item(140, 8).Value = "Sin especificar"
$ws.Cells.Item(140, 9).Value = "Primera"
$ws.Cells.Item(140, 10).Value = 100
$ws.Cells.Item(140, 11).Value = 20000
$ws.Cells.Item(140, 12).Value = 22000
$ws.Cells.Item(140, 13).Value = 21000
$ws.Cells.Item(140, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(140, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(140, 16).Value = 420
$ws.Cells.Item(140, 17).Value = 50
$ws.Cells.Item(140, 18).Value = "Hortaliza"
